# Add three new bulleted list items ("getAttributes() : Attributes",
# "getAttribute(Attribute)", "setAttribute(Attribute, Value)") right after
# the second "getOccurringContexts(S, P, O)" list entry (the one that is
# part of the numId=4 bullet list), matching its list/paragraph formatting.

$d = $word.ActiveDocument

# Locate the target paragraph: the "getOccurringContexts(S, P, O)" bullet
# that belongs to a numbered/bulleted list (ListType <> 0). The first
# "getOccurringContexts(S, P, O)" paragraph earlier in the document belongs
# to a different list, so keep scanning and remember the last match, which
# is the one immediately preceding the "FCA contexts? ..." paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd("`r", "`a")
    if ($text -eq "getOccurringContexts(S, P, O)" -and $p.Range.ListFormat.ListType -ne 0) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'getOccurringContexts(S, P, O)' list paragraph to anchor the new bullets on."
}

$newTexts = @(
    "getAttributes() : Attributes",
    "getAttribute(Attribute)",
    "setAttribute(Attribute, Value)"
)

$insertAfter = $targetIndex
foreach ($t in $newTexts) {
    $anchor = $d.Paragraphs.Item($insertAfter)
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($insertAfter + 1)
    $newPara.Range.Text = $t
    $insertAfter = $insertAfter + 1
}
